$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.197.27"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.825.64"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'234.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").Value = "'0.5997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.28%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "'0.07016"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.75%  "
$ws.Range("D9").Value = "'0.2774"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.19%  "
$ws.Range("D10").Value = "'23.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.72%  "
$ws.Range("D11").Value = "'0.07650"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "1.830.89"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "'4.779"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.98%  "
$ws.Range("D14").Value = "'0.000009890"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "'0.6224"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.99%  "
$ws.Range("D16").Value = "'78.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("D17").Value = "29.213.02"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "'5.814"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.82%  "
$ws.Range("D19").Value = "'222.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "'11.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.74%  "
$ws.Range("D22").Value = "'6.969"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.93%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").Value = "'155.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").Value = "'7.942"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.50%  "
$ws.Range("D26").Value = "'0.1289"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.65%  "
$ws.Range("D27").Value = "'16.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.09%  "
$ws.Range("D28").Value = "'1.478"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").Value = "'0.06214"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.97%  "
$ws.Range("D30").Value = "'1.437"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").Value = "'3.816"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.44%  "
$ws.Range("D32").Value = "'3.774"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.67%  "
$ws.Range("D33").Value = "'1.104"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("D34").Value = "'1.736"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("D35").Value = "'0.6419"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.17%  "
$ws.Range("D36").Value = "'2.541"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").Value = "1.218.81"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").Value = "'2.733"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.99%  "
$ws.Range("D39").Value = "'6.526"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.91%  "
$ws.Range("D40").Value = "'0.01723"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("D41").Value = "'0.8953"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.89%  "
$ws.Range("D42").Value = "'1.003"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "1.983.62"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "'100.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "'62.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.14%  "
$ws.Range("D46").Value = "'0.00000000115"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("D47").Value = "'8.510"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("D48").Value = "'0.4556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05494"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.564"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.35%  "
$ws.Range("D51").Value = "'6.380"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.50%  "
